$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.01786774324274964
$ws.Range("G2").Value = 38
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.01269327433982799
$ws.Range("G3").Value = 44
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.01120860421719229
$ws.Range("G4").Value = 40
$ws.Range("C5").Value = 88
$ws.Range("E5").Value = 0.01106086238608255
$ws.Range("G5").Value = 44
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.01161690410015671
$ws.Range("G6").Value = 44

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.01786774324274964
$ws.Range("G2").Value = 62
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.01269327433982799
$ws.Range("G3").Value = 62
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.01120860421719229
$ws.Range("G4").Value = 63
$ws.Range("C5").Value = 88
$ws.Range("E5").Value = 0.01106086238608255
$ws.Range("G5").Value = 64
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.01161690410015671
$ws.Range("G6").Value = 63

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.01786774324274964
$ws.Range("G2").Value = 69
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.01269327433982799
$ws.Range("G3").Value = 70
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.01120860421719229
$ws.Range("G4").Value = 70
$ws.Range("C5").Value = 88
$ws.Range("E5").Value = 0.01106086238608255
$ws.Range("G5").Value = 70
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.01161690410015671
$ws.Range("G6").Value = 70

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.01786774324274964
$ws.Range("G2").Value = 80
$ws.Range("C3").Value = 87
$ws.Range("E3").Value = 0.01269327433982799
$ws.Range("G3").Value = 81
$ws.Range("C4").Value = 87
$ws.Range("E4").Value = 0.01120860421719229
$ws.Range("G4").Value = 80
$ws.Range("C5").Value = 88
$ws.Range("E5").Value = 0.01106086238608255
$ws.Range("G5").Value = 80
$ws.Range("C6").Value = 88
$ws.Range("E6").Value = 0.01161690410015671
$ws.Range("G6").Value = 82

